# VIC_COVID19_DAILY_DATA.xlsx — append six new daily rows (109-114) to Sheet1,
# covering 2022-12-09 .. 2022-12-14 (serials 44912-44917), and move the
# frozen-pane/selection to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone formatting (number formats / styles) from the last existing
#        data row (108) down into the new block, so the new cells pick up the
#        same style indices (date format on col A, #,##0 on B/C/H) instead of
#        Excel minting brand-new style entries.
$ws.Range("A108:H108").Copy($ws.Range("A109:H114"))

# --- 2. Row 109: 2022-12-09
$ws.Range("A109").Value = 44912
$ws.Range("B109").Value = 2807
$ws.Range("D109").Value = 2055
$ws.Range("E109").Value = 623
$ws.Range("F109").Value = 26
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 5812

# --- 3. Row 110: 2022-12-10
$ws.Range("A110").Value = 44913
$ws.Range("B110").Value = 2989
$ws.Range("D110").Value = 2423
$ws.Range("E110").Value = 623
$ws.Range("F110").Value = 26
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 4509

# --- 4. Row 111: 2022-12-11
$ws.Range("A111").Value = 44914
$ws.Range("B111").Value = 3830
$ws.Range("D111").Value = 3109
$ws.Range("E111").Value = 668
$ws.Range("F111").Value = 33
$ws.Range("G111").Value = 20
$ws.Range("H111").Value = 5824

# --- 5. Row 112: 2022-12-12
$ws.Range("A112").Value = 44915
$ws.Range("B112").Value = 4176
$ws.Range("D112").Value = 3140
$ws.Range("E112").Value = 691
$ws.Range("F112").Value = 31
$ws.Range("G112").Value = 22
$ws.Range("H112").Value = 8097

# --- 6. Row 113: 2022-12-13
$ws.Range("A113").Value = 44916
$ws.Range("B113").Value = 3578
$ws.Range("D113").Value = 2580
$ws.Range("E113").Value = 704
$ws.Range("F113").Value = 30
$ws.Range("G113").Value = 22
$ws.Range("H113").Value = 7339

# --- 7. Row 114: 2022-12-14
$ws.Range("A114").Value = 44917
$ws.Range("B114").Value = 3243
$ws.Range("D114").Value = 2292
$ws.Range("E114").Value = 707
$ws.Range("F114").Value = 35
$ws.Range("G114").Value = 21
$ws.Range("H114").Value = 7290

# --- 8. Column C ("PCR cases") is a computed column (New Cases - In Hospital).
#        Fill it as one shared formula across the new rows, matching the
#        existing C101:C108 shared-formula pattern.
$ws.Range("C109:C114").Formula = "=B109-D109"

# --- 9. Move the frozen pane / active selection down to the new last row,
#        mirroring what Excel does automatically when new rows are appended
#        and the last row is re-selected. (Wrapped in [void] so the boolean
#        return value of Select() doesn't leak onto the output stream.)
[void]$ws.Range("A114").Select()
